# "Error Calculations and Plots"
#
# The sheet holds an ID column (A) plus four measured/derived numeric
# columns (B..F). This commit re-runs the error calculation used to decide
# which column-F values count as "present" vs "missing", and drops two
# samples (RM 232, SC 92) from the set entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (col 6) gets a fresh set of computed/missing values -------
# Newly computed (previously missing) values:
$ws.Range("F2").Value = 18.03
$ws.Range("F11").Value = 17.65
$ws.Range("F21").Value = 16.58

# Previously-present values that are now treated as missing:
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("F25").ClearContents()

# --- Two samples are removed outright, shifting everything below up -----
# Row 26 = "RM 232"
$ws.Rows("26:26").Delete()
# Row 28 ("SC 92") is now row 27 after the delete above
$ws.Rows("27:27").Delete()

# --- The last remaining sample ("SC 232") gets its D/E error + F value ---
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

Write-Output "missing_data.xlsx error calculations updated"
